$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C13 (row 14) placement values in place ---
$ws.Range("B14").Value = 127.7
$ws.Range("C14").Value = -112.3
$ws.Range("D14").Value = 0

# --- Insert a new row at 20 for new component "C19" (placed right after C18) ---
# This pushes D1-D4, J1,J2,J4, JP1, L1, Q1-Q3, R1-R12, U1-U4 down by one row.
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "C19"
$ws.Range("B20").Value = 128.12
$ws.Range("C20").Value = -110.18000000000001
$ws.Range("D20").Value = 180
$ws.Range("E20").Value = "top"

# --- Insert three new rows after R12 (now at row 44) for R13, R14, R15 ---
# This pushes U1-U4 (now at rows 45-48) further down by three rows.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = "R13"
$ws.Range("B45").Value = 129.50999999999999
$ws.Range("C45").Value = -107.3
$ws.Range("D45").Value = 90
$ws.Range("E45").Value = "top"

$ws.Range("A46").Value = "R14"
$ws.Range("B46").Value = 129.47999999999999
$ws.Range("C46").Value = -105.22
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = "top"

$ws.Range("A47").Value = "R15"
$ws.Range("B47").Value = 129.25999999999999
$ws.Range("C47").Value = -103.92
$ws.Range("D47").Value = -90
$ws.Range("E47").Value = "top"

# --- Remove U1 entirely (now at row 48); U2-U4 shift up to 48-50 ---
$ws.Rows.Item(48).Delete()

# --- Update U2 (now row 48) with its new placement values ---
$ws.Range("B48").Value = 125.15000000000001
$ws.Range("C48").Value = -105.72499999999999
